# Iceland Premier League - update match ordering/base data
# This script swaps/rotates the full data payload (columns B:AC) between the
# rows listed below. Column A (the sequential row id, 0-based) is left
# untouched on every row - only the underlying match data moves between rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple pairwise swaps: each pair exchanges its B:AC contents.
$swapPairs = @(
    @(2, 3),
    @(4, 6),
    @(34, 35),
    @(41, 42),
    @(67, 68),
    @(69, 71),
    @(88, 89)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = "B$r1" + ":AC$r1"
    $range2 = "B$r2" + ":AC$r2"

    $data1 = $ws.Range($range1).Value2
    $data2 = $ws.Range($range2).Value2

    $ws.Range($range1).Value2 = $data2
    $ws.Range($range2).Value2 = $data1
}

# Rotation groups: row N takes the B:AC content that used to live in the row
# referenced below (captured from the *original* sheet before any writes).
$rotationGroups = @(
    @(129, 130, 131),
    @(146, 147, 150, 149)
)

foreach ($group in $rotationGroups) {
    $count = $group.Length

    # Capture all of the original row contents first.
    $originals = @()
    foreach ($r in $group) {
        $rng = "B$r" + ":AC$r"
        $originals += , ($ws.Range($rng).Value2)
    }

    # new(group[i]) = old(group[i+1]), wrapping around to the start.
    for ($i = 0; $i -lt $count; $i++) {
        $destRow = $group[$i]
        $srcIndex = ($i + 1) % $count
        $destRange = "B$destRow" + ":AC$destRow"
        $ws.Range($destRange).Value2 = $originals[$srcIndex]
    }
}
